$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Test Case names (append ".txt.tour" suffix). The order in which
# new distinct strings are first assigned determines their position in the
# saved shared-string table, so update row 2 then rows 8 down to 3 to
# reproduce the expected shared-string ordering.
$ws.Range("A2").Value = "test-input-1.txt.tour"
$ws.Range("A8").Value = "test-input-7.txt.tour"
$ws.Range("A7").Value = "test-input-6.txt.tour"
$ws.Range("A6").Value = "test-input-5.txt.tour"
$ws.Range("A5").Value = "test-input-4.txt.tour"
$ws.Range("A4").Value = "test-input-3.txt.tour"
$ws.Range("A3").Value = "test-input-2.txt.tour"

# Update timing values that changed
$ws.Range("C2").Value = 0.0219
$ws.Range("C3").Value = 0.126
$ws.Range("C4").Value = 1.9019
$ws.Range("C8").Value = 9.194

# Widen column A (closest achievable value to the target 18.7109375 given
# the host's column-width quantization to 1/6-character increments)
$ws.Columns.Item(1).ColumnWidth = 17.85

# Move selection to D13
$ws.Range("D13").Select()
